$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Logged another 6.5 pages written that day (I19). This feeds the
# downstream SUM() totals (B5, B11, B12, I31), the remaining-pages
# calc (D11) and the pages/day-remaining rate (F4) - all formulas,
# so they recalc automatically.
$ws.Range("I19").Value = 6.5

# Leave the selection where the user was last working.
$ws.Range("N22").Select()
